$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.566.28'
$ws.Range("E2").Value = '  +0.55%  '

$ws.Range("D3").Value = '1.903.99'
$ws.Range("E3").Value = '  -0.64%  '

$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.31%  '

$ws.Range("D5").Value = "'337.90"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +4.01%  '

$ws.Range("E6").Value = '  -0.19%  '

$ws.Range("D7").Value = "'0.4768"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.93%  '

$ws.Range("D8").Value = "'0.4000"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.66%  '

$ws.Range("D9").Value = "'0.08054"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.14%  '

$ws.Range("D10").Value = "'0.9922"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.00%  '

$ws.Range("D11").Value = "'23.20"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.50%  '

$ws.Range("D12").Value = '1.884.29'
$ws.Range("E12").Value = '  -2.13%  '

$ws.Range("D13").Value = "'5.932"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.18%  '

$ws.Range("D14").Value = "'7.120"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.65%  '

$ws.Range("D15").Value = "'89.24"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.47%  '

$ws.Range("D16").Value = "'0.06814"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.33%  '

$ws.Range("E17").Value = '  -0.21%  '

$ws.Range("E18").Value = '  -1.44%  '

$ws.Range("D19").Value = "'17.35"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.28%  '

$ws.Range("D20").Value = "'1.005"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.11%  '

$ws.Range("D21").Value = '29.551.88'
$ws.Range("E21").Value = '  +0.46%  '

$ws.Range("D22").Value = "'5.507"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.51%  '

$ws.Range("D23").Value = "'11.68"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.46%  '

$ws.Range("D24").Value = "'2.158"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.86%  '

$ws.Range("D25").Value = '2.141.88'
$ws.Range("E25").Value = '  -0.88%  '

$ws.Range("D26").Value = "'156.78"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.82%  '

$ws.Range("D27").Value = "'6.537"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.72%  '

$ws.Range("D28").Value = "'19.61"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.94%  '

$ws.Range("D29").Value = "'2.059"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.54%  '

$ws.Range("D30").Value = "'119.29"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.00%  '

$ws.Range("D31").Value = "'0.9966"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.96%  '

$ws.Range("D32").Value = "'0.09544"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.59%  '

$ws.Range("D33").Value = "'5.469"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.23%  '

$ws.Range("D34").Value = "'1.389"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.25%  '

$ws.Range("D35").Value = "'3.532"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.49%  '

$ws.Range("D36").Value = "'0.06480"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +6.24%  '

$ws.Range("D37").Value = "'0.02243"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.64%  '

$ws.Range("D38").Value = "'1.196"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.54%  '

$ws.Range("D39").Value = "'0.5822"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.62%  '

$ws.Range("D40").Value = "'10.53"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.91%  '

$ws.Range("E41").Value = '  -4.23%  '

$ws.Range("D42").Value = "'0.1822"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.27%  '

$ws.Range("D43").Value = "'2.460"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.93%  '

$ws.Range("D44").Value = "'1.267"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.52%  '

$ws.Range("D45").Value = "'12.17"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.08%  '

$ws.Range("D46").Value = "'0.07401"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.53%  '

$ws.Range("D47").Value = "'0.5482"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.91%  '

$ws.Range("D48").Value = "'1.957"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.28%  '

$ws.Range("D49").Value = "'115.95"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.60%  '

$ws.Range("D50").Value = "'2.379"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.90%  '

$ws.Range("D51").Value = "'71.10"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.51%  '
